# Updated cryptos list - apply new Price (column D) and Volume(1h) (column E) values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a "Price" cell (column D) to a text value without letting Excel
# auto-convert numeric-looking strings (e.g. "596.85", "1.00") into real
# numbers, which would lose formatting (trailing zeros) or add float noise.
function Set-Price {
    param($row, $text)
    $cell = $ws.Range("D$row")
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Helper: set a "Volume(1h)" cell (column E). These are never numeric-looking
# (they include a leading/trailing double space and a percent sign) so a
# plain assignment keeps them as text automatically.
function Set-Volume {
    param($row, $pct)
    $ws.Range("E$row").Value = "  $pct  "
}

# Row 2 - Bitcoin
Set-Price  2 "65.746.36"
Set-Volume 2 "+0.14%"

# Row 3 - Ethereum
Set-Price  3 "2.664.70"
Set-Volume 3 "-0.19%"

# Row 4 - TetherUSD
Set-Price  4 "1.00"

# Row 5 - BNB
Set-Price  5 "596.85"
Set-Volume 5 "-0.84%"

# Row 6 - Solana
Set-Price  6 "157.07"
Set-Volume 6 "+0.00%"

# Row 7 - XRP
Set-Price  7 "0.659"
Set-Volume 7 "+6.74%"

# Row 8 - USDC
Set-Price  8 "1.00"
Set-Volume 8 "-0.01%"

# Row 9 - Dogecoin
Set-Volume 9 "-4.04%"

# Row 10 - Cardano
Set-Price  10 "0.400"
Set-Volume 10 "-0.63%"

# Row 11 - Toncoin
Set-Volume 11 "-0.39%"

# Row 12 - TRON
Set-Volume 12 "+1.36%"

# Row 13 - Avalanche
Set-Price  13 "28.78"
Set-Volume 13 "-2.21%"

# Row 14 - ShibaInu
Set-Price  14 "0.0000194"
Set-Volume 14 "-2.85%"

# Row 15 - Wrapped liquid staked Ether 2.0
Set-Price  15 "3.139.67"
Set-Volume 15 "-0.38%"

# Row 16 - WrappedBTC
Set-Price  16 "65.597.43"
Set-Volume 16 "+0.29%"

# Row 17 - WrappedEther
Set-Price  17 "2.659.50"
Set-Volume 17 "-0.99%"

# Row 18 - Chainlink
Set-Volume 18 "-2.50%"

# Row 19 - Polkadot
Set-Volume 19 "-0.64%"

# Row 20 - BitcoinCash
Set-Volume 20 "-0.71%"

# Row 21 - Uniswap
Set-Price  21 "7.42"
Set-Volume 21 "-3.54%"

# Row 22 - Dai
Set-Volume 22 "+0.05%"

# Row 23 - Litecoin
Set-Price  23 "69.57"
Set-Volume 23 "-0.37%"

# Row 24 - SuiNetwork
Set-Volume 24 "+8.85%"

# Row 25 - PEPE
Set-Volume 25 "+0.96%"

# Row 26 - InternetComputer(DFINITY)
Set-Price  26 "9.52"
Set-Volume 26 "-2.47%"

# Row 27 - Fetch.AI
Set-Volume 27 "+1.79%"

# Row 28 - Bittensor
Set-Price  28 "563.69"
Set-Volume 28 "+5.96%"

# Row 29 - Aptos
Set-Volume 29 "-0.81%"

# Row 30 - Kaspa
Set-Volume 30 "-2.83%"

# Row 31 - PancakeSwap
Set-Volume 31 "+0.25%"

# Row 32 - Binance-PegBSC-USD
Set-Price  32 "0.998"
Set-Volume 32 "-0.30%"

# Row 33 - ImmutableX
Set-Price  33 "1.80"
Set-Volume 33 "+2.70%"

# Row 34 - RenderToken
Set-Price  34 "6.55"
Set-Volume 34 "+0.23%"

# Row 35 - NEARProtocol
Set-Price  35 "5.47"
Set-Volume 35 "-0.79%"

# Row 36 - PolygonEcosystemToken
Set-Volume 36 "-0.89%"

# Row 37 - EthereumClassic
Set-Price  37 "20.46"
Set-Volume 37 "+0.06%"

# Row 38 - FirstDigitalUSD
Set-Volume 38 "-0.07%"

# Row 40 - Monero
Set-Price  40 "154.54"
Set-Volume 40 "-2.63%"

# Row 41 - Aave
Set-Price  41 "160.59"
Set-Volume 41 "-2.70%"

# Row 42 - Filecoin
Set-Volume 42 "-1.21%"

# Row 43 - Hedera
Set-Volume 43 "-1.13%"

# Row 44 - dogwifhat
Set-Volume 44 "-2.08%"

# Row 45 - InjectiveProtocol
Set-Price  45 "22.63"
Set-Volume 45 "-2.10%"

# Row 46 - Mantle
Set-Volume 46 "-0.80%"

# Row 47 - Stellar
Set-Price  47 "0.103"
Set-Volume 47 "+2.36%"

# Row 48 - VeChain
Set-Volume 48 "-1.92%"

# Row 49 - EnergySwap
Set-Price  49 "19.79"
Set-Volume 49 "-2.15%"

# Row 50 - BabyDogeCoin
Set-Volume 50 "+1.91%"

# Row 51 - ONDO
Set-Volume 51 "-2.06%"
